# Fruta / hortaliza, semanal
# Insert two new weekly price rows for Frutilla (Terminal Hortofrutícola Agro Chillán)
# at the top of the existing price history block, shifting the rest of the
# rows down by two (dimension grows from A1:T596 to A1:T598).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 494..596 down to 496..598, leaving rows 494:495 empty for the
# new records.
$ws.Rows("494:495").Insert()

# New row 494: Especial/Primera quality, bandeja 7 kilos, Provincia de Melipilla
$ws.Range("A494").Value = 7
$ws.Range("B494").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C494").Value = "Ñuble"
$ws.Range("D494").Value = 45218
$ws.Range("E494").Value = 16
$ws.Range("F494").Value = "Fruta"
$ws.Range("G494").Value = 100101
$ws.Range("H494").Value = "Berries"
$ws.Range("I494").Value = 100112025
$ws.Range("J494").Value = "Frutilla"
$ws.Range("K494").Value = "Sin especificar"
$ws.Range("L494").Value = "Primera"
$ws.Range("M494").Value = 60
$ws.Range("N494").Value = 10000
$ws.Range("O494").Value = 10000
$ws.Range("P494").Value = 10000
$ws.Range("Q494").Value = "$/bandeja 7 kilos"
$ws.Range("R494").Value = "Provincia de Melipilla"
$ws.Range("S494").Value = 1429
$ws.Range("T494").Value = 7

# New row 495: Segunda quality, bandeja 7 kilos, Provincia de Melipilla
$ws.Range("A495").Value = 7
$ws.Range("B495").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C495").Value = "Ñuble"
$ws.Range("D495").Value = 45218
$ws.Range("E495").Value = 16
$ws.Range("F495").Value = "Fruta"
$ws.Range("G495").Value = 100101
$ws.Range("H495").Value = "Berries"
$ws.Range("I495").Value = 100112025
$ws.Range("J495").Value = "Frutilla"
$ws.Range("K495").Value = "Sin especificar"
$ws.Range("L495").Value = "Segunda"
$ws.Range("M495").Value = 60
$ws.Range("N495").Value = 8000
$ws.Range("O495").Value = 8000
$ws.Range("P495").Value = 8000
$ws.Range("Q495").Value = "$/bandeja 7 kilos"
$ws.Range("R495").Value = "Provincia de Melipilla"
$ws.Range("S495").Value = 1143
$ws.Range("T495").Value = 7
